$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Hora") bumps from 10 to 11 for every data row (2-51).
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Range("G$r")
    $cell.Value = "'11"
    $cell.Style = "Normal"
}

# Per-cell updates to Coin / Link / Price / Volume(1h) columns.
$ws.Range("D2").Value = "'330.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.20%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.642"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.93%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08209"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.10%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = "'8.761"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = "'2.012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.28%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.85%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.18%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9218"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.00%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1274"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.34%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1949"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09407"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.63%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03844"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'7.96%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1060"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.02%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006110"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.02%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D21").Value = "'8.301"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-4.93%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1374"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.63%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.2662"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'6.31%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04407"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.72%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001257"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.004315"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-6.27%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-2.37%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02755"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'9.72%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05511"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.28%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007900"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'5.66%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.13%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008950"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-9.72%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'2.65%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01143"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.31%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006780"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.55%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003192"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'5.14%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.01%"
$ws.Range("E51").Style = "Normal"
